# Added repeatability to Battle Interrupts
# Insert a new repeatable "check-register-style" interrupt test block
# right after the existing "check-register-root" interrupt block
# (before the old row 19 / END_SCENE), shifting subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 19, pushing everything
# from the old row 19 onward down to row 22 onward.
$ws.Rows.Item(19).Resize(3).Insert()

# New INTERRUPT block (rows 19-21)
$ws.Cells.Item(19, 1).Value = "INTERRUPT"
$ws.Cells.Item(19, 2).Value = "check-register-style"

$ws.Cells.Item(20, 1).Value = "DIALOGUE"
$ws.Cells.Item(20, 2).Value = "Clarke"
$ws.Cells.Item(20, 3).Value = "You registered the root {last-cast,style}. Stylish!"
$ws.Cells.Item(20, 4).Value = "_"
$ws.Cells.Item(20, 5).Value = "tanuki"
$ws.Cells.Item(20, 6).Value = "CENTER"
$ws.Cells.Item(20, 7).Value = "END_DIALOGUE"

$ws.Cells.Item(21, 1).Value = "END_INTERRUPT"

# Update the selection to reflect where the editor left off (C20),
# matching the post-edit sheetView selection in the diff.
$ws.Range("C20").Select()
